# Lecture 22 handout - slide 15 edits
#  - trace-table text corrections
#  - reflow / reposition several highlight rectangles
#  - add a new highlight rectangle (duplicate of "Rectangle 25") with its
#    own click-triggered "Appear" animation

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)

# ---------------------------------------------------------------------
# 1. Table text edits (shape "Table 3")
# ---------------------------------------------------------------------
$tbl = $s.Shapes.Item("Table 3").Table

# Row 7 / Col 1 : "(after 2nd DONE, RET)->R6 " -> "(after 2nd DONE, before RET)->R6 "
$tbl.Cell(7,1).Shape.TextFrame.TextRange.Text = "(after 2nd DONE, before RET)->R6 "

# Row 7 / Col 3 : "R3= [R1+2]=[6004]= 0 (NULL)" -> "R3= [R1+2]=[6005]= 0 (NULL)"
$tbl.Cell(7,3).Shape.TextFrame.TextRange.Text = "R3= [R1+2]=[6005]= 0 (NULL)"

# Row 11 / Col 3 : empty cell -> " R6 (before 2nd RET)"
$tbl.Cell(11,3).Shape.TextFrame.TextRange.Text = " R6 (before 2nd RET)"

# ---------------------------------------------------------------------
# 2. Reposition / resize existing highlight rectangles
# ---------------------------------------------------------------------

# Rectangle 15 : off (840700,2669874)/ext (2669315,355516) -> off (533400,2669874)/ext (2976615,335029)
$r15 = $s.Shapes.Item("Rectangle 15")
$r15.Left   = 42.0000393701
$r15.Top    = 210.2263385827
$r15.Width  = 234.3791732283
$r15.Height = 26.3802755906

# Rectangle 21 : off (6655957,4934682) -> (840700,4602745)
$r21 = $s.Shapes.Item("Rectangle 21")
$r21.Left = 66.1968897638
$r21.Top  = 362.4209055118

# Rectangle 22 : off (840699,4567351) -> (6655956,4230636)
$r22 = $s.Shapes.Item("Rectangle 22")
$r22.Left = 524.0910629921
$r22.Top  = 333.1209842520

# Rectangle 23 : off (3748329,6068575) -> (6655957,4975454)
$r23 = $s.Shapes.Item("Rectangle 23")
$r23.Left = 524.0911417323
$r23.Top  = 391.7680708661

# Rectangle 24 : off (840700,5056466) -> (3748329,6075657)
$r24 = $s.Shapes.Item("Rectangle 24")
$r24.Left = 295.1440551181
$r24.Top  = 478.3982283465

# Rectangle 25 : off (840698,5367823) -> (840700,4958261)
$r25 = $s.Shapes.Item("Rectangle 25")
$r25.Left = 66.1968897638
$r25.Top  = 390.4142913386

# ---------------------------------------------------------------------
# 3. Add a new highlight rectangle ("Rectangle 26"), a duplicate of
#    "Rectangle 25", placed where "Rectangle 25" used to sit.
# ---------------------------------------------------------------------

# burn through the low, currently-unused shape ids (2,3,13,20) so the
# duplicate below lands on id 27, matching "Rectangle 26"'s real id
for ($i = 1; $i -le 4; $i++) {
    $dummy = $s.Shapes.AddShape(1, 0, 0, 1, 1)
    $dummy.Delete()
}

$dupRange = $r25.Duplicate()
$r26 = $dupRange.Item(1)
$r26.Name   = "Rectangle 26"
$r26.Left   = 66.1968897638
$r26.Top    = 421.0812992126
$r26.Width  = 210.1823228346
$r26.Height = 27.9934251969

# ---------------------------------------------------------------------
# 4. Animate the new rectangle the same way as its siblings: an
#    on-click "Appear" entrance effect.
# ---------------------------------------------------------------------
$seq = $s.TimeLine.MainSequence
$null = $seq.AddEffect($r26, 1)
